$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting D:K data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from column F into the new D:E columns
$ws.Range("F1").EntireColumn.Copy()
$ws.Range("D1:E1").EntireColumn.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = Dec-2018, E = Sep-2018)
$ws.Range("D7").Value = 43462
$ws.Range("E7").Value = 43371
$ws.Range("D8").Value = 4265300
$ws.Range("E8").Value = 3913600
$ws.Range("D9").Value = 3794400
$ws.Range("E9").Value = 3383800
$ws.Range("D10").Value = 470900
$ws.Range("E10").Value = 529800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = -157300
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 150700
$ws.Range("E15").Value = 152500
$ws.Range("D17").Value = 3892000
$ws.Range("E17").Value = 3631100
$ws.Range("D18").Value = 373300
$ws.Range("E18").Value = 282500
$ws.Range("D20").Value = -82900
$ws.Range("E20").Value = -92500
$ws.Range("D21").Value = 441100
$ws.Range("E21").Value = 342400
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 290400
$ws.Range("E23").Value = 189900
$ws.Range("D24").Value = 51000
$ws.Range("E24").Value = 68300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 239400
$ws.Range("E26").Value = 121600
$ws.Range("D27").Value = 239400
$ws.Range("E27").Value = 121500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 11300
$ws.Range("E29").Value = 54000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 82900
$ws.Range("E32").Value = 92500
$ws.Range("D33").Value = 250700
$ws.Range("E33").Value = 175500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 250700
$ws.Range("E35").Value = 175500
$ws.Range("D38").Value = 43462
$ws.Range("E38").Value = 43371
$ws.Range("D41").Value = 249900
$ws.Range("E41").Value = 215000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1880300
$ws.Range("E43").Value = 1790400
$ws.Range("D44").Value = 371100
$ws.Range("E44").Value = 724800
$ws.Range("D45").Value = 148700
$ws.Range("E45").Value = 171200
$ws.Range("D46").Value = 2650000
$ws.Range("E46").Value = 2901400
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2153200
$ws.Range("E48").Value = 1378100
$ws.Range("D49").Value = 7605500
$ws.Range("E49").Value = 7747400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1330300
$ws.Range("E52").Value = 1693200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 13738900
$ws.Range("E54").Value = 13720100
$ws.Range("D57").Value = 866200
$ws.Range("E57").Value = 1018900
$ws.Range("D58").Value = 53400
$ws.Range("E58").Value = 30900
$ws.Range("D59").Value = 1277700
$ws.Range("E59").Value = 1440300
$ws.Range("D60").Value = 2197300
$ws.Range("E60").Value = 2490200
$ws.Range("D61").Value = 7323700
$ws.Range("E61").Value = 7213100
$ws.Range("D62").Value = 1000100
$ws.Range("E62").Value = 987300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 10521000
$ws.Range("E66").Value = 10690500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 990400
$ws.Range("E72").Value = 710500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3217900
$ws.Range("E76").Value = 3029600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43462
$ws.Range("E80").Value = 43371
$ws.Range("D81").Value = 250700
$ws.Range("E81").Value = 175500
$ws.Range("D83").Value = 150700
$ws.Range("E83").Value = 152500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -207400
$ws.Range("E89").Value = 905900
$ws.Range("D91").Value = -114400
$ws.Range("E91").Value = -195800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 194200
$ws.Range("E94").Value = -193100
$ws.Range("D96").Value = -27200
$ws.Range("E96").Value = -25800
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 48100
$ws.Range("E100").Value = -663800
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 34900
$ws.Range("E102").Value = 49100

# Row 14 (Minority Interest) - newly reported quarters are NA, not 0
$ws.Range("F14:J14").Value = "NA"
